# Daily refresh of the cryptos price/volume table (Price column D,
# Volume(1h) column E). Values that look numeric (e.g. "32.91", "1.00")
# are written with a leading apostrophe so Excel keeps them as Text,
# matching the existing (and unchanged) cell formatting for this column
# instead of silently reformatting them as numbers (which would drop
# trailing zeros, e.g. turn "1.00" into "1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.663.58'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '3.688.84'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'667.65"
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('D6').Value = "'160.26"
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +3.45%  '
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('D12').Value = "'0.0000233"
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = "'32.91"
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = '3.666.28'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '69.656.61'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').Value = "'16.17"
$ws.Range('E17').Value = '  +0.88%  '
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').Value = "'470.35"
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D23').Value = '3.834.57'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').Value = "'0.0000127"
$ws.Range('E24').Value = '  +4.69%  '
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').Value = "'9.05"
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').Value = "'1.70"
$ws.Range('E29').Value = '  -2.87%  '
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = "'0.165"
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').Value = '3.680.15'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('D37').Value = "'6.12"
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('D39').Value = "'2.25"
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = "'177.01"
$ws.Range('E41').Value = '  +1.99%  '
$ws.Range('D42').Value = "'0.0907"
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').Value = "'47.02"
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').Value = "'2.76"
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('D46').Value = "'27.61"
$ws.Range('E46').Value = '  -2.14%  '
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E48').Value = '  -2.00%  '
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('E51').Value = '  -0.25%  '
